# Actualizacion automatica del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update results for already-played matches (rows 14-16): mark as "Fallo" / -1 ---
$ws.Range("G14").Value = "Fallo"
$ws.Range("H14").Value = -1

$ws.Range("G15").Value = "Fallo"
$ws.Range("H15").Value = -1

$ws.Range("G16").Value = "Fallo"
$ws.Range("H16").Value = -1

# --- Append new pending matches (rows 18-35) ---
$newRows = @(
    @(14639119, "2025-09-08", "Alycia Parks", "Darja Vidmanova", "Gana Alycia Parks", 1.67),
    @(14633471, "2025-09-08", "Carolina Alves", "Nauhany Vitoria Leme Da Silva", "Gana Nauhany Vitoria Leme Da Silva", 4),
    @(14633474, "2025-09-08", "Julia Riera", "Vitalia Diatchenko", "Gana Julia Riera", 1.73),
    @(14633455, "2025-09-08", "Dimitar Kuzmanov", "Alan Wazny", "Gana Alan Wazny", 10),
    @(14633448, "2025-09-09", "Jan Choinski", "Vitaliy Sachko", "Gana Vitaliy Sachko", 2.38),
    @(14633447, "2025-09-09", "Tiago Pereira", "Andrew Paulson", "Gana Andrew Paulson", 2),
    @(14633326, "2025-09-09", "Antoine Escoffier", "Patrick Zahraj", "Gana Patrick Zahraj", 1.83),
    @(14633191, "2025-09-09", "Gregoire Barrere", "Daniel Jade", "Gana Daniel Jade", 8),
    @(14633440, "2025-09-09", "Filippo Romano", "Giovanni Fonio", "Gana Filippo Romano", 2.5),
    @(14633446, "2025-09-09", "Stefano Napolitano", "Enrico Dalla Valle", "Gana Stefano Napolitano", 1.73),
    @(14631273, "2025-09-09", "Daniel Evans", "Charles Chen", "Gana Charles Chen", 8),
    @(14633144, "2025-09-08", "Franco Agamenone", "Max Alcala Gurri", "Gana Max Alcala Gurri", 1.83),
    @(14633141, "2025-09-09", "Sebastian Gima", "Andrea Picchione", "Gana Sebastian Gima", 2.1),
    @(14634120, "2025-09-08", "Blaise Bicknell", "Johannus Monday", "Gana Blaise Bicknell", 3),
    @(14633836, "2025-09-08", "Emon van Loben Sels", "Martin Damm Jr", "Gana Emon van Loben Sels", 3.5),
    @(14634112, "2025-09-08", "Andre Ilagan", "Trevor Svajda", "Gana Trevor Svajda", 1.57),
    @(14634118, "2025-09-08", "Patrick Maloney", "Luca Pow", "Gana Luca Pow", 3),
    @(14637407, "2025-09-08", "Lorenzo Carboni", "Iannis Miletich", "Gana Iannis Miletich", 3)
)

$row = 18
foreach ($data in $newRows) {
    $ws.Cells.Item($row, 1).Value = $data[0]
    # Leading apostrophe forces the date to be stored as plain text,
    # matching the "fecha" text formatting used throughout column B.
    $ws.Cells.Item($row, 2).Value = "'" + $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
    $ws.Cells.Item($row, 5).Value = $data[4]
    $ws.Cells.Item($row, 6).Value = $data[5]
    $row++
}
